$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Araraquara covid.csv")

# ------------------------------------------------------------------
# 1. Preserve the "latest row" highlight formatting (style of row 503,
#    columns C and E:L) by copying it onto row 510 BEFORE row 503 is
#    restyled to the normal (non-highlighted) look.
# ------------------------------------------------------------------
$ws.Range("C503").Copy() | Out-Null
$ws.Range("C510").PasteSpecial(-4122) | Out-Null
$ws.Range("E503:L503").Copy() | Out-Null
$ws.Range("E510:L510").PasteSpecial(-4122) | Out-Null

# ------------------------------------------------------------------
# 2. Row 503 is no longer the last row with data, so it goes back to
#    the regular (non-highlighted) formatting used by row 502. Also
#    fix the K503 typo (6943 -> 5943).
# ------------------------------------------------------------------
$ws.Range("C502").Copy() | Out-Null
$ws.Range("C503").PasteSpecial(-4122) | Out-Null
$ws.Range("E502:L502").Copy() | Out-Null
$ws.Range("E503:L503").PasteSpecial(-4122) | Out-Null
$ws.Range("K503").Value = 5943

# ------------------------------------------------------------------
# 3. Apply the regular row formatting to the newly used rows 504-509
#    (row 510 already received the highlight formatting above).
# ------------------------------------------------------------------
$ws.Range("C502").Copy() | Out-Null
$ws.Range("C504:C509").PasteSpecial(-4122) | Out-Null
$ws.Range("E502:L502").Copy() | Out-Null
$ws.Range("E504:L504").PasteSpecial(-4122) | Out-Null
$ws.Range("E505:L505").PasteSpecial(-4122) | Out-Null
$ws.Range("E506:L506").PasteSpecial(-4122) | Out-Null
$ws.Range("E507:L507").PasteSpecial(-4122) | Out-Null
$ws.Range("E508:L508").PasteSpecial(-4122) | Out-Null
$ws.Range("E509:L509").PasteSpecial(-4122) | Out-Null

# ------------------------------------------------------------------
# 4. Fill in the new daily figures for 15-21 Aug 2021 (rows 504-510).
#    Column D keeps the "novos" (new cases) shared formula; the rest
#    are plain data values pulled from the daily bulletin.
# ------------------------------------------------------------------
$ws.Range("C504").Value = 29027
$ws.Range("E504").Value = 568
$ws.Range("F504").Value = 52
$ws.Range("G504").Value = 51
$ws.Range("H504").Value = 28
$ws.Range("I504").Value = 165247
$ws.Range("J504").Value = 64890
$ws.Range("K504").Value = 5943
$ws.Range("L504").Value = 236080

$ws.Range("C505").Value = 29034
$ws.Range("E505").Value = 569
$ws.Range("F505").Value = 54
$ws.Range("G505").Value = 51
$ws.Range("H505").Value = 27
$ws.Range("I505").Value = 165300
$ws.Range("J505").Value = 64950
$ws.Range("K505").Value = 5943
$ws.Range("L505").Value = 236193

$ws.Range("C506").Value = 29114
$ws.Range("E506").Value = 569
$ws.Range("F506").Value = 47
$ws.Range("G506").Value = 46
$ws.Range("H506").Value = 26
$ws.Range("I506").Value = 167382
$ws.Range("J506").Value = 66139
$ws.Range("K506").Value = 5944
$ws.Range("L506").Value = 239465

$ws.Range("C507").Value = 29167
$ws.Range("E507").Value = 569
$ws.Range("F507").Value = 48
$ws.Range("G507").Value = 47
$ws.Range("H507").Value = 29
$ws.Range("I507").Value = 169292
$ws.Range("J507").Value = 67944
$ws.Range("K507").Value = 5944
$ws.Range("L507").Value = 243180

$ws.Range("C508").Value = 29194
$ws.Range("E508").Value = 569
$ws.Range("F508").Value = 43
$ws.Range("G508").Value = 42
$ws.Range("H508").Value = 24
$ws.Range("I508").Value = 171294
$ws.Range("J508").Value = 69946
$ws.Range("K508").Value = 5944
$ws.Range("L508").Value = 247184

$ws.Range("C509").Value = 29225
$ws.Range("E509").Value = 570
$ws.Range("F509").Value = 38
$ws.Range("G509").Value = 38
$ws.Range("H509").Value = 22
$ws.Range("I509").Value = 173042
$ws.Range("J509").Value = 71889
$ws.Range("K509").Value = 5944
$ws.Range("L509").Value = 250875

$ws.Range("C510").Value = 29271
$ws.Range("E510").Value = 570
$ws.Range("F510").Value = 41
$ws.Range("G510").Value = 41
$ws.Range("H510").Value = 23
$ws.Range("I510").Value = 174035
$ws.Range("J510").Value = 74818
$ws.Range("K510").Value = 5977
$ws.Range("L510").Value = 254830

# Extend the "novos" shared formula (D column, = Cn - Cn-1) down to row 510.
$ws.Range("D504:D510").Formula = "=(C504-C503)"
